# Auto-generated script updating TPM-derived NATMI metrics for Vegfa-Kdr sheet
# Applies the new TPM-based values from the commit "update scripts wuth new tpm"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 2.162809
$ws.Range("H2").Value = 6.488427000000001
$ws.Range("I2").Value = 0.06755089002018773
$ws.Range("J2").Value = 0.06755089002018773
$ws.Range("M2").Value = 181.3526613333333
$ws.Range("N2").Value = 544.057984
$ws.Range("O2").Value = 0.9845849379007657
$ws.Range("P2").Value = 0.984584937900766
$ws.Range("Q2").Value = 392.2311681056854
$ws.Range("R2").Value = 3530.080512951169
$ws.Range("S2").Value = 0.06650958885566799
$ws.Range("T2").Value = 0.066509588855668
# Row 3
$ws.Range("G3").Value = 2.162809
$ws.Range("H3").Value = 6.488427000000001
$ws.Range("I3").Value = 0.06755089002018773
$ws.Range("J3").Value = 0.06755089002018773
$ws.Range("O3").Value = 0.003686045149950483
$ws.Range("P3").Value = 0.003686045149950484
$ws.Range("Q3").Value = 1.46841754246
$ws.Range("R3").Value = 13.21575788214
$ws.Range("S3").Value = 0.0002489956305337515
$ws.Range("T3").Value = 0.0002489956305337515
# Row 4
$ws.Range("G4").Value = 2.162809
$ws.Range("H4").Value = 6.488427000000001
$ws.Range("I4").Value = 0.06755089002018773
$ws.Range("J4").Value = 0.06755089002018773
$ws.Range("M4").Value = 0.6398506666666667
$ws.Range("N4").Value = 1.919552
$ws.Range("O4").Value = 0.003473824559694892
$ws.Range("P4").Value = 0.003473824559694892
$ws.Range("Q4").Value = 1.383874780522667
$ws.Range("R4").Value = 12.454873024704
$ws.Range("S4").Value = 0.0002346599407813767
$ws.Range("T4").Value = 0.0002346599407813767
# Row 5
$ws.Range("G5").Value = 2.162809
$ws.Range("H5").Value = 6.488427000000001
$ws.Range("I5").Value = 0.06755089002018773
$ws.Range("J5").Value = 0.06755089002018773
$ws.Range("M5").Value = 1.520540333333333
$ws.Range("N5").Value = 4.561621
$ws.Range("O5").Value = 0.008255192389588805
$ws.Range("P5").Value = 0.008255192389588807
$ws.Range("Q5").Value = 3.288638317796333
$ws.Range("R5").Value = 29.597744860167
$ws.Range("S5").Value = 0.0005576455932046041
$ws.Range("T5").Value = 0.0005576455932046042
# Row 6
$ws.Range("I6").Value = 0.5628021396814664
$ws.Range("J6").Value = 0.5628021396814664
$ws.Range("M6").Value = 181.3526613333333
$ws.Range("N6").Value = 544.057984
$ws.Range("O6").Value = 0.9845849379007657
$ws.Range("P6").Value = 0.984584937900766
$ws.Range("Q6").Value = 3267.885006306646
$ws.Range("R6").Value = 29410.96505675981
$ws.Range("S6").Value = 0.5541265097486947
$ws.Range("T6").Value = 0.5541265097486948
# Row 7
$ws.Range("I7").Value = 0.5628021396814664
$ws.Range("J7").Value = 0.5628021396814664
$ws.Range("O7").Value = 0.003686045149950483
$ws.Range("P7").Value = 0.003686045149950484
$ws.Range("S7").Value = 0.002074514097354624
$ws.Range("T7").Value = 0.002074514097354624
# Row 8
$ws.Range("I8").Value = 0.5628021396814664
$ws.Range("J8").Value = 0.5628021396814664
$ws.Range("M8").Value = 0.6398506666666667
$ws.Range("N8").Value = 1.919552
$ws.Range("O8").Value = 0.003473824559694892
$ws.Range("P8").Value = 0.003473824559694892
$ws.Range("Q8").Value = 11.52979164740267
$ws.Range("R8").Value = 103.768124826624
$ws.Range("S8").Value = 0.001955075895074313
$ws.Range("T8").Value = 0.001955075895074313
# Row 9
$ws.Range("I9").Value = 0.5628021396814664
$ws.Range("J9").Value = 0.5628021396814664
$ws.Range("M9").Value = 1.520540333333333
$ws.Range("N9").Value = 4.561621
$ws.Range("O9").Value = 0.008255192389588805
$ws.Range("P9").Value = 0.008255192389588807
$ws.Range("Q9").Value = 27.39938261866133
$ws.Range("R9").Value = 246.594443567952
$ws.Range("S9").Value = 0.004646039940342737
$ws.Range("T9").Value = 0.004646039940342738
# Row 10
$ws.Range("G10").Value = 4.650307000000001
$ws.Range("H10").Value = 13.950921
$ws.Range("I10").Value = 0.1452427730405732
$ws.Range("J10").Value = 0.1452427730405732
$ws.Range("M10").Value = 181.3526613333333
$ws.Range("N10").Value = 544.057984
$ws.Range("O10").Value = 0.9845849379007657
$ws.Range("P10").Value = 0.984584937900766
$ws.Range("Q10").Value = 843.3455504670295
$ws.Range("R10").Value = 7590.109954203265
$ws.Range("S10").Value = 0.1430038466746878
$ws.Range("T10").Value = 0.1430038466746879
# Row 11
$ws.Range("G11").Value = 4.650307000000001
$ws.Range("H11").Value = 13.950921
$ws.Range("I11").Value = 0.1452427730405732
$ws.Range("J11").Value = 0.1452427730405732
$ws.Range("O11").Value = 0.003686045149950483
$ws.Range("P11").Value = 0.003686045149950484
$ws.Range("Q11").Value = 3.15727943458
$ws.Range("R11").Value = 28.41551491122
$ws.Range("S11").Value = 0.0005353714191315637
$ws.Range("T11").Value = 0.0005353714191315639
# Row 12
$ws.Range("G12").Value = 4.650307000000001
$ws.Range("H12").Value = 13.950921
$ws.Range("I12").Value = 0.1452427730405732
$ws.Range("J12").Value = 0.1452427730405732
$ws.Range("M12").Value = 0.6398506666666667
$ws.Range("N12").Value = 1.919552
$ws.Range("O12").Value = 0.003473824559694892
$ws.Range("P12").Value = 0.003473824559694892
$ws.Range("Q12").Value = 2.975502034154667
$ws.Range("R12").Value = 26.779518307392
$ws.Range("S12").Value = 0.0005045479121065344
$ws.Range("T12").Value = 0.0005045479121065344
# Row 13
$ws.Range("G13").Value = 4.650307000000001
$ws.Range("H13").Value = 13.950921
$ws.Range("I13").Value = 0.1452427730405732
$ws.Range("J13").Value = 0.1452427730405732
$ws.Range("M13").Value = 1.520540333333333
$ws.Range("N13").Value = 4.561621
$ws.Range("O13").Value = 0.008255192389588805
$ws.Range("P13").Value = 0.008255192389588807
$ws.Range("Q13").Value = 7.070979355882334
$ws.Range("R13").Value = 63.638814202941
$ws.Range("S13").Value = 0.001199007034647314
$ws.Range("T13").Value = 0.001199007034647315
# Row 14
$ws.Range("G14").Value = 7.184856000000001
$ws.Range("H14").Value = 21.554568
$ws.Range("I14").Value = 0.2244041972577726
$ws.Range("J14").Value = 0.2244041972577726
$ws.Range("M14").Value = 181.3526613333333
$ws.Range("N14").Value = 544.057984
$ws.Range("O14").Value = 0.9845849379007657
$ws.Range("P14").Value = 0.984584937900766
$ws.Range("Q14").Value = 1302.992756896768
$ws.Range("R14").Value = 11726.93481207091
$ws.Range("S14").Value = 0.2209449926217152
$ws.Range("T14").Value = 0.2209449926217153
# Row 15
$ws.Range("G15").Value = 7.184856000000001
$ws.Range("H15").Value = 21.554568
$ws.Range("I15").Value = 0.2244041972577726
$ws.Range("J15").Value = 0.2244041972577726
$ws.Range("O15").Value = 0.003686045149950483
$ws.Range("P15").Value = 0.003686045149950484
$ws.Range("Q15").Value = 4.87808613264
$ws.Range("R15").Value = 43.90277519376001
$ws.Range("S15").Value = 0.0008271640029305443
$ws.Range("T15").Value = 0.0008271640029305444
# Row 16
$ws.Range("G16").Value = 7.184856000000001
$ws.Range("H16").Value = 21.554568
$ws.Range("I16").Value = 0.2244041972577726
$ws.Range("J16").Value = 0.2244041972577726
$ws.Range("M16").Value = 0.6398506666666667
$ws.Range("N16").Value = 1.919552
$ws.Range("O16").Value = 0.003473824559694892
$ws.Range("P16").Value = 0.003473824559694892
$ws.Range("Q16").Value = 4.597234901504001
$ws.Range("R16").Value = 41.37511411353601
$ws.Range("S16").Value = 0.0007795408117326676
$ws.Range("T16").Value = 0.0007795408117326677
# Row 17
$ws.Range("G17").Value = 7.184856000000001
$ws.Range("H17").Value = 21.554568
$ws.Range("I17").Value = 0.2244041972577726
$ws.Range("J17").Value = 0.2244041972577726
$ws.Range("M17").Value = 1.520540333333333
$ws.Range("N17").Value = 4.561621
$ws.Range("O17").Value = 0.008255192389588805
$ws.Range("P17").Value = 0.008255192389588807
$ws.Range("Q17").Value = 10.924863337192
$ws.Range("R17").Value = 98.32377003472801
$ws.Range("S17").Value = 0.00185249982139415
$ws.Range("T17").Value = 0.00185249982139415

Write-Output "Updated 174 cells with new TPM values"
